$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("STAGE")
$ws2 = $wb.Worksheets.Item("LMSPROD")

# --- STAGE sheet (sheet1) new values ---
$ws1.Range("A2").Value = "FPK12School71871"
$ws1.Range("B2").Value = "FPK12Classroom77158"
$ws1.Range("C2").Value = "FPK12Section13169"
$ws1.Range("E3").Value = "96441"
$ws1.Range("E4").Value = "35549"
$ws1.Range("E5").Value = "77223"

# --- LMSPROD sheet (sheet2) new values ---
$ws2.Range("A2").Value = "FPK12School22646"
$ws2.Range("B2").Value = "FPK12Classroom56973"
$ws2.Range("C2").Value = "FPK12Section78958"
$ws2.Range("E3").Value = "36393"
$ws2.Range("E4").Value = "42259"
$ws2.Range("E5").Value = "7643"

# --- Selection / active sheet ---
$ws2.Range("E4").Select()
$ws2.Activate()
